$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.565.87"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.590.07"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "609.15"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "147.89"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.488"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "8.04"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "4.202.58"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "29.99"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "3.597.51"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "66.643.64"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "11.49"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "428.46"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").Value = "79.01"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "3.739.12"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "0.0000121"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "8.24"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "9.29"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "3.588.80"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "25.47"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.45"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").Value = "0.157"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "1.71"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("D38").Value = "5.63"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "177.12"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").Value = "0.0858"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "0.897"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").Value = "2.56"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "25.01"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "24.21"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.950"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "2.414.05"
$ws.Range("E51").Value = "  +4.54%  "
